# Updated cryptos list on Sat Aug 31 17:00:58 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "58.849.70"
$ws.Range("E2").Value = "  +0.77%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.493.86"
$ws.Range("E3").Value = "  +1.66%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
Set-TextValue "D5" "533.75"
$ws.Range("E5").Value = "  +1.04%  "

# Row 6 - Solana
Set-TextValue "D6" "136.04"
$ws.Range("E6").Value = "  +1.59%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.10%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.63%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.506.47"
$ws.Range("E9").Value = "  +1.86%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.30%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.76%  "

# Row 12 - Toncoin
Set-TextValue "D12" "5.38"
$ws.Range("E12").Value = "  +1.42%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +1.48%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "2.942.13"
$ws.Range("E14").Value = "  +1.91%  "

# Row 15 - now WrappedBTC (was Avalanche)
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D15" "58.714.58"
$ws.Range("E15").Value = "  +0.69%  "

# Row 16 - now Avalanche (was WrappedBTC)
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D16" "22.79"
$ws.Range("E16").Value = "  +0.87%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -0.50%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.503.24"
$ws.Range("E18").Value = "  +1.78%  "

# Row 19 - Chainlink
Set-TextValue "D19" "11.03"
$ws.Range("E19").Value = "  +2.90%  "

# Row 20 - Polkadot
Set-TextValue "D20" "4.24"
$ws.Range("E20").Value = "  +1.18%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "322.41"
$ws.Range("E21").Value = "  +0.53%  "

# Row 23 - Uniswap
Set-TextValue "D23" "5.96"
$ws.Range("E23").Value = "  +4.10%  "

# Row 24 - Litecoin
Set-TextValue "D24" "65.25"
$ws.Range("E24").Value = "  +4.37%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  +3.13%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -0.25%  "

# Row 27 - Binance-PegBSC-USD
Set-TextValue "D27" "0.995"
$ws.Range("E27").Value = "  +1.31%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "7.49"
$ws.Range("E28").Value = "  +0.71%  "

# Row 29 - PEPE
Set-TextValue "D29" "0.0₃0761"
$ws.Range("E29").Value = "  +1.71%  "

# Row 30 - Aptos
Set-TextValue "D30" "6.52"
$ws.Range("E30").Value = "  +0.77%  "

# Row 31 - Monero
Set-TextValue "D31" "171.33"
$ws.Range("E31").Value = "  +4.68%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -0.22%  "

# Row 33 - Fetch.AI
Set-TextValue "D33" "1.16"
$ws.Range("E33").Value = "  +7.15%  "

# Row 34 - USDe
Set-TextValue "D34" "0.998"
$ws.Range("E34").Value = "  +0.03%  "

# Row 35 - EthereumClassic
Set-TextValue "D35" "18.29"
$ws.Range("E35").Value = "  +0.44%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +0.07%  "

# Row 37 - NEARProtocol
Set-TextValue "D37" "4.03"
$ws.Range("E37").Value = "  +0.17%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  -0.80%  "

# Row 39 - OKB
Set-TextValue "D39" "36.75"
$ws.Range("E39").Value = "  +1.01%  "

# Row 40 - SuiNetwork
Set-TextValue "D40" "0.803"
$ws.Range("E40").Value = "  +0.91%  "

# Row 41 - Filecoin
Set-TextValue "D41" "3.56"
$ws.Range("E41").Value = "  +0.75%  "

# Row 42 - Bittensor
Set-TextValue "D42" "281.81"
$ws.Range("E42").Value = "  +2.38%  "

# Row 43 - RenderToken
Set-TextValue "D43" "5.21"
$ws.Range("E43").Value = "  +3.19%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  -0.11%  "

# Row 45 - Mantle
Set-TextValue "D45" "0.605"
$ws.Range("E45").Value = "  +3.40%  "

# Row 46 - Aave
Set-TextValue "D46" "129.79"
$ws.Range("E46").Value = "  +7.72%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -0.03%  "

# Row 49 - Hedera
Set-TextValue "D49" "0.0501"
$ws.Range("E49").Value = "  -0.93%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  +0.21%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "17.22"
$ws.Range("E51").Value = "  +0.90%  "
